# Slide 20, shape "Rectangle 3" (the second code-listing box) contains the
# line:
#     if (!matchTypes(variable.type, expr))
# as a paragraph built from several runs. The final run of that paragraph
# is ", expr))". The edit splits that single run into two runs, ", " and
# "expr))", leaving the visible text unchanged.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)
$shape = $s.Shapes.Item("Rectangle 3")
$tr = $shape.TextFrame.TextRange

$fullText = $tr.Text
$splitAt0 = $fullText.IndexOf("expr))")
if ($splitAt0 -lt 0) {
    throw "edit.ps1: could not locate 'expr))' in Rectangle 3 text"
}

# Characters() is 1-based.
$splitAt1 = $splitAt0 + 1
$tail = $tr.Characters($splitAt1, 6)

# Re-assigning the (identical) text of just this sub-range forces
# PowerPoint to break the previously-single run into two runs right at
# this boundary, while each new run keeps the original run formatting
# (Consolas / size 1800 / etc.).
$tail.Text = "expr))"
